# Update need_to_buy.xlsx values on the active worksheet (Sheet 1)
# per the R-generated refresh of columns B, C, E, F for rows 2-7 and 9-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10283.7123862307
$ws.Range("C2").Value = 9788.160728946
$ws.Range("E2").Value = 5504.89258898721
$ws.Range("F2").Value = -49.9311117527832

$ws.Range("B3").Value = 10479.163620429
$ws.Range("C3").Value = 9781.48275652265
$ws.Range("E3").Value = 5308.95811672708
$ws.Range("F3").Value = 253.626703052072

$ws.Range("B4").Value = 10503.9340317843
$ws.Range("C4").Value = 9393.80916483286
$ws.Range("E4").Value = 6105.01539279898
$ws.Range("F4").Value = 270.642689901327

$ws.Range("B5").Value = 4432.57332848585
$ws.Range("C5").Value = 7281.55637732575
$ws.Range("E5").Value = 6010.81699347466
$ws.Range("F5").Value = 178.70722378335

$ws.Range("B6").Value = 4516.05831396367
$ws.Range("C6").Value = 7591.96302386868
$ws.Range("E6").Value = 6283.2987301338
$ws.Range("F6").Value = 202.994239750103

$ws.Range("B7").Value = 12370.2746011533
$ws.Range("C7").Value = 10852.2345108451
$ws.Range("E7").Value = 7378.58433692522
$ws.Range("F7").Value = 384.475785323764

$ws.Range("C9").Value = 10491.1248498996
$ws.Range("F9").Value = 329.040024610629

$ws.Range("C10").Value = 10783.8563640652
$ws.Range("F10").Value = 341.237171034195

$ws.Range("C11").Value = 10223.2192252566
$ws.Range("F11").Value = 317.877290250504

$ws.Range("C12").Value = 7414.35093150552
$ws.Range("F12").Value = 185.165262463042

$ws.Range("C13").Value = 7375.38016896097
$ws.Range("F13").Value = 183.187789775277

$ws.Range("C14").Value = 11306.1544830264
$ws.Range("F14").Value = 422.093069679041

$ws.Range("C15").Value = 11344.1441649174
$ws.Range("F15").Value = 423.675973091169
